$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells
$ws.Range("B11").Value = 20363
$ws.Range("B12").Value = 20777

# Data for rows 13 through 23 (A..I)
$data = @(
    @(12, 20382, 0, 0, 12615, 0, 0, 0, 0),
    @(13, 21520, 0, 0, 13127, 0, 0, 0, 0),
    @(14, 22336, 0, 0, 13457, 0, 0, 0, 0),
    @(15, 22942, 0, 0, 14118, 0, 0, 0, 0),
    @(16, 23137, 0, 0, 14239, 0, 0, 0, 0),
    @(17, 22951, 0, 0, 14529, 0, 0, 0, 0),
    @(18, 24926, 0, 0, 15330, 0, 0, 0, 0),
    @(19, 25705, 0, 0, 16529, 0, 0, 0, 0),
    @(20, 24384, 0, 0, 15922, 0, 0, 0, 0),
    @(21, 24443, 0, 0, 16382, 0, 0, 0, 0),
    @(22, 24616, 0, 0, 16531, 0, 0, 0, 0)
)

$startRow = 13
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
